$d = $word.ActiveDocument

# --- Change 1: remove the standalone "Beta Version" paragraph, keeping the
# bookmarkStart/bookmarkEnd markers by merging them into the preceding
# (title) paragraph. ---
$r = $d.Content
$found = $r.Find.Execute("Beta Version", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Beta Version' paragraph" }
$betaPara = $r.Paragraphs(1)
$titlePara = $betaPara.Previous
$mergeRange = $d.Range($titlePara.Range.Start, $betaPara.Range.End)
$mergeRange.InsertXML('<w:p w14:paraId="6BE039F1" w14:textId="574E1A09" w:rsidR="004D6C85" w:rsidRDefault="004D6C85" w:rsidP="004D6C85"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00B1217F"><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Сопроводительная документация к приложению, реализующему функции необходимые в лабораторной работе </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# --- Change 2: split the image-format list so "jpeg" is its own run with
# English language formatting, and drop "tif" from the list. ---
$r = $d.Content
$found = $r.Find.Execute("jpg, gif, tif, bmp, png, pcx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find image-format run" }
$formatsPara = $r.Paragraphs(1)
$formatsPara.Range.InsertXML('<w:p w14:paraId="35F36AD0" w14:textId="70F23822" w:rsidR="004D6C85" w:rsidRDefault="004D6C85" w:rsidP="004D6C85"><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">В данном предложении представлены возможности работы с такими форматами изображений как </w:t></w:r><w:r w:rsidRPr="004D6C85"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">jpg, gif, bmp, png, </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>jpeg</w:t></w:r></w:p>')

# --- Change 3: append a sentence about where the temp folder for archive
# extraction is created. ---
$r = $d.Content
$found = $r.Find.Execute("Затем папка и все ее содержимое удаляется.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find archive-folder paragraph" }
$archivePara = $r.Paragraphs(1)
$archivePara.Range.InsertXML('<w:p w14:paraId="410C4085" w14:textId="77777777" w:rsidR="000274DC" w:rsidRDefault="004D6C85" w:rsidP="004D6C85"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Выбор изображений осуществляется с помощью функции </w:t></w:r><w:r w:rsidRPr="004D6C85"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>getOpenFileNames</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> класса</w:t></w:r><w:r w:rsidRPr="004D6C85"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>QFileDialog</w:t></w:r><w:r w:rsidRPr="004D6C85"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.Выбрать</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> можно либо сразу изображения либо архив формата </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>zip</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> .Работа с архивами представлена классами </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>qzipreader</w:t></w:r><w:r w:rsidRPr="004D6C85"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">и </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>qzipwriter</w:t></w:r><w:r w:rsidR="00982310" w:rsidRPr="00982310"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00982310"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">которые не являются частью </w:t></w:r><w:r w:rsidR="00982310"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Qt</w:t></w:r><w:r w:rsidR="00982310" w:rsidRPr="00982310"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00982310"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>API</w:t></w:r><w:r w:rsidR="00982310"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">, но могут быть подключены.  Чтобы получить информацию из архива он предварительно проходит проверку, есть ли хоть 1 файл нужного формата. При наличие такового создается временна папка в нее распаковываются все файлы архива и выбираются те, что подходят выбранному формату. Затем папка и все ее содержимое удаляется. </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Папка создается там же, где и находится архив.</w:t></w:r></w:p>')

# --- Change 4: append a note about the sorting limitation to the table
# paragraph. ---
$r = $d.Content
$found = $r.Find.Execute("а сортировка с помощью встроенных функций", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find table-sorting paragraph" }
$tablePara = $r.Paragraphs(1)
$tablePara.Range.InsertXML('<w:p w14:paraId="3566B590" w14:textId="207843BC" w:rsidR="004D6C85" w:rsidRPr="000274DC" w:rsidRDefault="000274DC" w:rsidP="004D6C85"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Таблица реализована в виде </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>QTableWidget</w:t></w:r><w:r w:rsidR="003C56E0" w:rsidRPr="003C56E0"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="003C56E0"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>а сортировка с помощью встроенных функций</w:t></w:r><w:r w:rsidRPr="000274DC"><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">В связи с некорректной работой сортировки строк, состоящих из чисел, сортировка для разрешения не работает, однако для размера реализована своя, которая работает исправно. </w:t></w:r></w:p>')
